$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Adolfo Jimenez Dimas, period 1805, valor mora 31249
$ws.Range("C16").Value = "73209005"
$ws.Range("D16").Value = "ADOLFO JIMENEZ DIMAS"
$ws.Range("E16").Value = "1805"
$ws.Range("F16").Value = 31249

# Row 17: Adolfo Jimenez Dimas, period 1804, valor mora 1042
$ws.Range("C17").Value = "73209005"
$ws.Range("D17").Value = "ADOLFO JIMENEZ DIMAS"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 1042

# Row 18: Carlos Alberto Cervantes Julio, period 1805, valor mora 31249
$ws.Range("C18").Value = "1047377965"
$ws.Range("D18").Value = "CARLOS ALBERTO CERVANTES JULIO"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 31249

# Row 19: Carlos Alberto Cervantes Julio, period 1804, valor mora 1042
$ws.Range("C19").Value = "1047377965"
$ws.Range("D19").Value = "CARLOS ALBERTO CERVANTES JULIO"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 1042
